$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: Date cell "26 June 2025" -> "20 February 2026"
#   Original runs: "2" | "6 " | "June" | " " | "2025"
#   Target runs:   "20" | " February" | " " | "2026"
#   (the last two runs - a space with spacing=-2, and the year with
#    spacing=-4 - are reused unmodified other than the year's text).
# ---------------------------------------------------------------------------
$dateParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7E0F76C8" w14:textId="46439975" w:rsidR="0028517B" w:rsidRDefault="00000000"><w:pPr><w:pStyle w:val="TableParagraph"/><w:spacing w:line="253" w:lineRule="exact"/><w:ind w:left="109"/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>20</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> February</w:t></w:r><w:r><w:rPr><w:spacing w:val="-2"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:spacing w:val="-4"/><w:sz w:val="24"/></w:rPr><w:t>2026</w:t></w:r></w:p>'

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "26 June 2025*") {
        $p.Range.InsertXML($dateParaXml)
        break
    }
}

# ---------------------------------------------------------------------------
# Change 2: Team ID "LTVIP2025TMID41465" -> "LTVIP2026TMIDS47801"
#   Single run, simple text substitution.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("LTVIP2025TMID41465", $true, $false, $false, $false, $false, $true, 1, $false, "LTVIP2026TMIDS47801", 2)

# ---------------------------------------------------------------------------
# Change 3: "Sorting:Tranfer" -> "Sorting :Transfer"
#   The mis-spelled run (wrapped in spellcheck/grammar proofErr markers) is
#   replaced by five plain runs with identical run formatting, and the
#   proofErr markers are dropped.
# ---------------------------------------------------------------------------
$sortParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4D9526C6" w14:textId="77777777" w:rsidR="0028517B" w:rsidRDefault="00000000"><w:pPr><w:pStyle w:val="TableParagraph"/><w:ind w:left="109"/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Smart</w:t></w:r><w:r><w:rPr><w:spacing w:val="-1"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Sorting</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>:Tran</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>fer</w:t></w:r><w:r><w:rPr><w:spacing w:val="-4"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Learning</w:t></w:r><w:r><w:rPr><w:spacing w:val="-1"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>for</w:t></w:r><w:r><w:rPr><w:spacing w:val="-4"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:spacing w:val="-2"/><w:sz w:val="24"/></w:rPr><w:t>Identifying</w:t></w:r></w:p>'

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Smart Sorting:Tranfer Learning for Identifying*") {
        $p.Range.InsertXML($sortParaXml)
        break
    }
}
